$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the old second header row (units row): this shifts the data rows
# (Chancy-Pougny / Verbois / Seujet) up so they sit directly under row 1.
$ws.Range("A2:K2").EntireRow.Delete()

# Wipe out whatever is left of the old top header row (values + styles)
# so we can rebuild it from scratch.
$ws.Range("A1:K1").ClearContents()
$ws.Range("A1:K1").ClearFormats()

# New header row.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# F1:K1 carry the small (9pt) Arial font used elsewhere in the header/labels.
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

$ws.Range("A2:K2").Select()
